$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 1114, shifting existing data down.
$ws.Rows.Item(1114).Insert()

$ws.Range("A1114").Value = 7
$ws.Range("B1114").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C1114").Value = "Ñuble"
$ws.Range("D1114").Value = 45021
$ws.Range("E1114").Value = 16
$ws.Range("F1114").Value = "Fruta"
$ws.Range("G1114").Value = 100104
$ws.Range("H1114").Value = "Frutos de pepita"
$ws.Range("I1114").Value = 100104002
$ws.Range("J1114").Value = "Manzana"
$ws.Range("K1114").Value = "Granny Smith"
$ws.Range("L1114").Value = "Especial"
$ws.Range("M1114").Value = 100
$ws.Range("N1114").Value = 12000
$ws.Range("O1114").Value = 12000
$ws.Range("P1114").Value = 12000
$ws.Range("Q1114").Value = "$/caja 16 kilos empedrada"
$ws.Range("R1114").Value = "Región de O'Higgins"
$ws.Range("S1114").Value = 750
$ws.Range("T1114").Value = 16
